# Ontbrekende velden zib Patient toegevoegd
# Adds the two missing zib Patient field mappings as new rows at the
# bottom of column B on the "map" sheet, and leaves the selection on
# the last newly-added cell (matching the saved worksheet state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Patient.GenderIdentity"
$ws.Range("B17").Value = "Patient.MultipleBirthSequence"

$ws.Range("B17").Select()
